$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2210.9736
$ws.Range("I15").Value = 2210.9736
$ws.Range("K15").Value = 6632.9208
$ws.Range("M15").Value = -6463.9208
$ws.Range("H28").Value = 1229.8889
$ws.Range("I28").Value = 1473.5
$ws.Range("J28").Value = 742.6667
$ws.Range("K28").Value = 1473.5
$ws.Range("L28").Value = 742.6667
$ws.Range("M28").Value = -988.5
$ws.Range("N28").Value = -1712.6667
$ws.Range("H40").Value = 2093.375
$ws.Range("I40").Value = 1498.3334
$ws.Range("J40").Value = 2230.6924
$ws.Range("K40").Value = 1498.3334
$ws.Range("L40").Value = 2230.6924
$ws.Range("M40").Value = -1323.3334
$ws.Range("N40").Value = -2580.6924
$ws.Range("H116").Value = 18984.125
$ws.Range("J116").Value = 18148.334
$ws.Range("L116").Value = 18148.334
$ws.Range("N116").Value = -25032.334
$ws.Range("H127").Value = 2349.4167
$ws.Range("I127").Value = 1999.25
$ws.Range("J127").Value = 2524.5
$ws.Range("K127").Value = 5997.75
$ws.Range("L127").Value = 7573.5
$ws.Range("M127").Value = -1037.75
$ws.Range("N127").Value = -17493.5
$ws.Range("H132").Value = 4016.9348
$ws.Range("I132").Value = 3879.282
$ws.Range("K132").Value = 11637.846
$ws.Range("M132").Value = -9107.846000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1214.6923
$ws.Range("I2").Value = 1093.5555
$ws.Range("K2").Value = 1093.5555
$ws.Range("M2").Value = -980.5554999999999
$ws.Range("H31").Value = 60000
$ws.Range("J31").Value = 60000
$ws.Range("L31").Value = 60000
$ws.Range("N31").Value = -60588
$ws.Range("H32").Value = 2861.48
$ws.Range("I32").Value = 2457.4443
$ws.Range("K32").Value = 2457.4443
$ws.Range("M32").Value = -2170.4443
$ws.Range("H45").Value = 1913.2354
$ws.Range("I45").Value = 1354
$ws.Range("K45").Value = 1354
$ws.Range("M45").Value = -977
$ws.Range("H88").Value = 4902.4
$ws.Range("I88").Value = 1118.7778
$ws.Range("J88").Value = 6523.952
$ws.Range("K88").Value = 1118.7778
$ws.Range("L88").Value = 6523.952
$ws.Range("M88").Value = -712.7778000000001
$ws.Range("N88").Value = -7335.952
$ws.Range("H91").Value = 4902.4
$ws.Range("I91").Value = 1118.7778
$ws.Range("J91").Value = 6523.952
$ws.Range("K91").Value = 1118.7778
$ws.Range("L91").Value = 6523.952
$ws.Range("M91").Value = 285.2221999999999
$ws.Range("N91").Value = -9331.952000000001
$ws.Range("H97").Value = 25887
$ws.Range("I97").Value = 100010
$ws.Range("J97").Value = 7356.25
$ws.Range("K97").Value = 100010
$ws.Range("L97").Value = 7356.25
$ws.Range("M97").Value = -99514
$ws.Range("N97").Value = -8348.25
$ws.Range("H103").Value = 73499.75
$ws.Range("J103").Value = 73499.75
$ws.Range("L103").Value = 73499.75
$ws.Range("N103").Value = -75843.75
$ws.Range("H116").Value = 1214.6923
$ws.Range("I116").Value = 1093.5555
$ws.Range("K116").Value = 1093.5555
$ws.Range("M116").Value = 1200.4445
$ws.Range("H132").Value = 3097.0312
$ws.Range("J132").Value = 7673.6665
$ws.Range("L132").Value = 23020.9995
$ws.Range("N132").Value = -28080.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1214.6923
$ws.Range("I3").Value = 1093.5555
$ws.Range("K3").Value = 1093.5555
$ws.Range("M3").Value = -979.5554999999999
$ws.Range("H20").Value = 1656.6154
$ws.Range("I20").Value = 1361.579
$ws.Range("J20").Value = 2457.4285
$ws.Range("K20").Value = 1361.579
$ws.Range("L20").Value = 2457.4285
$ws.Range("M20").Value = -1114.579
$ws.Range("N20").Value = -2951.4285
$ws.Range("H94").Value = 17704.46
$ws.Range("I94").Value = 1724.5
$ws.Range("J94").Value = 24806.666
$ws.Range("K94").Value = 1724.5
$ws.Range("L94").Value = 24806.666
$ws.Range("M94").Value = -1273.5
$ws.Range("N94").Value = -25708.666
$ws.Range("H107").Value = 1678.4651
$ws.Range("I107").Value = 1744.9445
$ws.Range("J107").Value = 1336.5714
$ws.Range("K107").Value = 1744.9445
$ws.Range("L107").Value = 1336.5714
$ws.Range("M107").Value = 175.0554999999999
$ws.Range("N107").Value = -5176.5714
$ws.Range("H115").Value = 85500
$ws.Range("J115").Value = 85500
$ws.Range("L115").Value = 85500
$ws.Range("N115").Value = -88634

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2436666.8
$ws.Range("I6").Value = 2436666.8
$ws.Range("K6").Value = 2436666.8
$ws.Range("M6").Value = -2436553.8
$ws.Range("H31").Value = 3062.4167
$ws.Range("I31").Value = 1144.3334
$ws.Range("J31").Value = 4213.2666
$ws.Range("K31").Value = 1144.3334
$ws.Range("L31").Value = 4213.2666
$ws.Range("M31").Value = -849.3334
$ws.Range("N31").Value = -4803.2666
$ws.Range("H34").Value = 3062.4167
$ws.Range("I34").Value = 1144.3334
$ws.Range("J34").Value = 4213.2666
$ws.Range("K34").Value = 1144.3334
$ws.Range("L34").Value = 4213.2666
$ws.Range("M34").Value = -942.3334
$ws.Range("N34").Value = -4617.2666
$ws.Range("H99").Value = 3065.3809
$ws.Range("I99").Value = 2959.2144
$ws.Range("K99").Value = 2959.2144
$ws.Range("M99").Value = -1461.2144
$ws.Range("H107").Value = 14302.75
$ws.Range("J107").Value = 2494.5
$ws.Range("L107").Value = 2494.5
$ws.Range("N107").Value = -6334.5
$ws.Range("H122").Value = 5114.278
$ws.Range("I122").Value = 5074.643
$ws.Range("J122").Value = 5253
$ws.Range("K122").Value = 15223.929
$ws.Range("L122").Value = 15759
$ws.Range("M122").Value = -12773.929
$ws.Range("N122").Value = -20659
$ws.Range("H126").Value = 3065.3809
$ws.Range("I126").Value = 2959.2144
$ws.Range("K126").Value = 8877.643199999999
$ws.Range("M126").Value = -6407.643199999999
$ws.Range("H134").Value = 5900.449
$ws.Range("I134").Value = 5667.273
$ws.Range("J134").Value = 6381.375
$ws.Range("K134").Value = 17001.819
$ws.Range("L134").Value = 19144.125
$ws.Range("M134").Value = -14466.819
$ws.Range("N134").Value = -24214.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 101204330
$ws.Range("I4").Value = 112338140
$ws.Range("K4").Value = 337014420
$ws.Range("M4").Value = -337014308
$ws.Range("H129").Value = 11907274
$ws.Range("I129").Value = 1326.3334
$ws.Range("J129").Value = 20836734
$ws.Range("K129").Value = 3979.0002
$ws.Range("L129").Value = 62510202
$ws.Range("M129").Value = 1020.9998
$ws.Range("N129").Value = -62520202
$ws.Range("H131").Value = 1125922.9
$ws.Range("J131").Value = 3650.5762
$ws.Range("L131").Value = 10951.7286
$ws.Range("N131").Value = -21031.7286
$ws.Range("H132").Value = 4334.6665
$ws.Range("J132").Value = 4334.6665
$ws.Range("L132").Value = 39011.9985
$ws.Range("N132").Value = -44071.9985
$ws.Range("H134").Value = 11084.56
$ws.Range("I134").Value = 5162.615
$ws.Range("J134").Value = 17500
$ws.Range("K134").Value = 15487.845
$ws.Range("L134").Value = 52500
$ws.Range("M134").Value = -10417.845
$ws.Range("N134").Value = -62640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 13800
$ws.Range("H80").Value = 2067.9092
$ws.Range("I80").Value = 2058
$ws.Range("J80").Value = 2079.8
$ws.Range("K80").Value = 2058
$ws.Range("L80").Value = 2079.8
$ws.Range("M80").Value = -1060
$ws.Range("N80").Value = -4075.8
$ws.Range("H83").Value = 2067.9092
$ws.Range("I83").Value = 2058
$ws.Range("J83").Value = 2079.8
$ws.Range("K83").Value = 10290
$ws.Range("L83").Value = 10399
$ws.Range("M83").Value = -5298
$ws.Range("N83").Value = -20383
$ws.Range("H122").Value = 2635.125
$ws.Range("I122").Value = 2233
$ws.Range("K122").Value = 6699
$ws.Range("M122").Value = -4249

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6155.3335
$ws.Range("I132").Value = 6000
$ws.Range("J132").Value = 6186.4
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 18559.2
$ws.Range("M132").Value = -15470
$ws.Range("N132").Value = -23619.2
$ws.Range("H136").Value = 1497.1666
$ws.Range("I136").Value = 1497.1666
$ws.Range("K136").Value = 4491.4998
$ws.Range("M136").Value = -1941.4998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4187.4287
$ws.Range("I81").Value = 4228.174
$ws.Range("J81").Value = 4000
$ws.Range("K81").Value = 8456.348
$ws.Range("L81").Value = 8000
$ws.Range("M81").Value = -7395.348
$ws.Range("N81").Value = -10122
$ws.Range("H84").Value = 4187.4287
$ws.Range("I84").Value = 4228.174
$ws.Range("J84").Value = 4000
$ws.Range("K84").Value = 42281.74
$ws.Range("L84").Value = 40000
$ws.Range("M84").Value = -36977.74
$ws.Range("N84").Value = -50608
$ws.Range("H132").Value = 9930.444
$ws.Range("I132").Value = 16099.25
$ws.Range("J132").Value = 4995.4
$ws.Range("K132").Value = 48297.75
$ws.Range("L132").Value = 14986.2
$ws.Range("M132").Value = -45767.75
$ws.Range("N132").Value = -20046.2
